# Fruta / hortaliza, semanal
#
# A new weekly record is inserted as row 42 (pushing the existing rows
# 42-87 down to 43-88). The new record is a near-duplicate of the record
# that used to sit at row 42 (same market/region/product/quality), but
# with a newer date and updated prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 42, shifting rows
# 42:87 down to 43:88 (dimension grows from A1:R87 to A1:R88).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = 'Macroferia Regional de Talca'
$ws.Range("C42").Value = 'Maule'
$ws.Range("D42").Value = 44894
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 300000000
$ws.Range("G42").Value = 'Espárragos'
$ws.Range("H42").Value = 'Sin especificar'
$ws.Range("I42").Value = 'Primera'
$ws.Range("J42").Value = 3000
$ws.Range("K42").Value = 1100
$ws.Range("L42").Value = 1100
$ws.Range("M42").Value = 1100
$ws.Range("N42").Value = '$/kilo'
$ws.Range("O42").Value = 'Provincia de Linares'
$ws.Range("P42").Value = 1100
$ws.Range("Q42").Value = 1
$ws.Range("R42").Value = 'Hortaliza'
